# Auto-generated edit script: apply scheduled-runner market-data refresh
# to the Ragnarok_Profits workbook. Updates currentAveragePrice /
# LevePrice / LeveProfit columns (H-N) for affected leve rows across
# the ALC, ARM, BSM, CRP, CUL, LTW and WVR crafting-job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1055.8948
$ws.Range("I4").Value = 566.9
$ws.Range("K4").Value = 566.9
$ws.Range("M4").Value = -452.9
$ws.Range("H9").Value = 302.4
$ws.Range("I9").Value = 212.5
$ws.Range("J9").Value = 437.25
$ws.Range("K9").Value = 212.5
$ws.Range("L9").Value = 437.25
$ws.Range("M9").Value = -43.5
$ws.Range("N9").Value = -775.25
$ws.Range("H15").Value = 987.81134
$ws.Range("I15").Value = 987.81134
$ws.Range("K15").Value = 2963.43402
$ws.Range("M15").Value = -2794.43402
$ws.Range("H17").Value = 1628.5555
$ws.Range("J17").Value = 1632.8
$ws.Range("L17").Value = 4898.4
$ws.Range("N17").Value = -5234.4
$ws.Range("H19").Value = 1295.6
$ws.Range("I19").Value = 1038.4286
$ws.Range("J19").Value = 1520.625
$ws.Range("K19").Value = 1038.4286
$ws.Range("L19").Value = 1520.625
$ws.Range("M19").Value = -863.4286
$ws.Range("N19").Value = -1870.625
$ws.Range("H33").Value = 169.9
$ws.Range("I33").Value = 174.88889
$ws.Range("K33").Value = 174.88889
$ws.Range("M33").Value = 54.11111
$ws.Range("H40").Value = 55571932
$ws.Range("J40").Value = 100028380
$ws.Range("L40").Value = 100028380
$ws.Range("N40").Value = -100028730
$ws.Range("H63").Value = 99999
$ws.Range("J63").Value = 99999
$ws.Range("L63").Value = 99999
$ws.Range("N63").Value = -101247
$ws.Range("H64").Value = 7028.9287
$ws.Range("J64").Value = 7198.1816
$ws.Range("L64").Value = 7198.1816
$ws.Range("N64").Value = -7694.1816
$ws.Range("H66").Value = 99999
$ws.Range("J66").Value = 99999
$ws.Range("L66").Value = 299997
$ws.Range("N66").Value = -306237
$ws.Range("H67").Value = 7028.9287
$ws.Range("J67").Value = 7198.1816
$ws.Range("L67").Value = 7198.1816
$ws.Range("N67").Value = -8914.1816
$ws.Range("H74").Value = 3957.4
$ws.Range("I74").Value = 3957.4
$ws.Range("K74").Value = 3957.4
$ws.Range("M74").Value = -3021.4
$ws.Range("H77").Value = 3957.4
$ws.Range("I77").Value = 3957.4
$ws.Range("K77").Value = 19787
$ws.Range("M77").Value = -15107
$ws.Range("H103").Value = 38462980
$ws.Range("I103").Value = 837.8
$ws.Range("K103").Value = 2513.4
$ws.Range("M103").Value = -1927.4
$ws.Range("H132").Value = 1282.0698
$ws.Range("I132").Value = 1198.25
$ws.Range("J132").Value = 2399.6667
$ws.Range("K132").Value = 3594.75
$ws.Range("L132").Value = 7199.000100000001
$ws.Range("M132").Value = -1064.75
$ws.Range("N132").Value = -12259.0001
$ws.Range("H135").Value = 8518.857
$ws.Range("I135").Value = 265.66666
$ws.Range("K135").Value = 2390.99994
$ws.Range("M135").Value = 144.0000600000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7728.451
$ws.Range("I32").Value = 6410.366
$ws.Range("J32").Value = 16236.091
$ws.Range("K32").Value = 6410.366
$ws.Range("L32").Value = 16236.091
$ws.Range("M32").Value = -6123.366
$ws.Range("N32").Value = -16810.091
$ws.Range("H74").Value = 1384.8182
$ws.Range("I74").Value = 1337.1
$ws.Range("J74").Value = 1862
$ws.Range("K74").Value = 1337.1
$ws.Range("L74").Value = 1862
$ws.Range("M74").Value = -463.0999999999999
$ws.Range("N74").Value = -3610
$ws.Range("H77").Value = 1384.8182
$ws.Range("I77").Value = 1337.1
$ws.Range("J77").Value = 1862
$ws.Range("K77").Value = 6685.5
$ws.Range("L77").Value = 9310
$ws.Range("M77").Value = -2317.5
$ws.Range("N77").Value = -18046
$ws.Range("H97").Value = 2581.3572
$ws.Range("I97").Value = 1729
$ws.Range("J97").Value = 3717.8333
$ws.Range("K97").Value = 1729
$ws.Range("L97").Value = 3717.8333
$ws.Range("M97").Value = -1233
$ws.Range("N97").Value = -4709.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1369.7333
$ws.Range("J80").Value = 1475.0834
$ws.Range("L80").Value = 1475.0834
$ws.Range("N80").Value = -3471.0834
$ws.Range("H83").Value = 1369.7333
$ws.Range("J83").Value = 1475.0834
$ws.Range("L83").Value = 7375.416999999999
$ws.Range("N83").Value = -17359.417
$ws.Range("H86").Value = 2278.48
$ws.Range("I86").Value = 1582.45
$ws.Range("K86").Value = 1582.45
$ws.Range("M86").Value = -459.45
$ws.Range("H89").Value = 2278.48
$ws.Range("I89").Value = 1582.45
$ws.Range("K89").Value = 7912.25
$ws.Range("M89").Value = -2296.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 77.045456
$ws.Range("I7").Value = 175.5
$ws.Range("J7").Value = 55.166668
$ws.Range("K7").Value = 175.5
$ws.Range("L7").Value = 55.166668
$ws.Range("M7").Value = -62.5
$ws.Range("N7").Value = -281.166668
$ws.Range("H22").Value = 327.75
$ws.Range("I22").Value = 187
$ws.Range("K22").Value = 187
$ws.Range("M22").Value = 163
$ws.Range("H31").Value = 21279318
$ws.Range("I31").Value = 31252266
$ws.Range("J31").Value = 3695.8
$ws.Range("K31").Value = 31252266
$ws.Range("L31").Value = 3695.8
$ws.Range("M31").Value = -31251971
$ws.Range("N31").Value = -4285.8
$ws.Range("H34").Value = 21279318
$ws.Range("I34").Value = 31252266
$ws.Range("J34").Value = 3695.8
$ws.Range("K34").Value = 31252266
$ws.Range("L34").Value = 3695.8
$ws.Range("M34").Value = -31252064
$ws.Range("N34").Value = -4099.8
$ws.Range("H99").Value = 14224.15
$ws.Range("I99").Value = 7519.8184
$ws.Range("K99").Value = 7519.8184
$ws.Range("M99").Value = -6021.8184
$ws.Range("H126").Value = 14224.15
$ws.Range("I126").Value = 7519.8184
$ws.Range("K126").Value = 22559.4552
$ws.Range("M126").Value = -20089.4552
$ws.Range("H132").Value = 2516
$ws.Range("I132").Value = 2237.8125
$ws.Range("K132").Value = 6713.4375
$ws.Range("M132").Value = -4183.4375
$ws.Range("H134").Value = 1705.2703
$ws.Range("J134").Value = 2723.889
$ws.Range("L134").Value = 8171.667
$ws.Range("N134").Value = -13241.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 147391520
$ws.Range("I4").Value = 170290110
$ws.Range("K4").Value = 510870330
$ws.Range("M4").Value = -510870218
$ws.Range("H42").Value = 166673950
$ws.Range("J42").Value = 10328.25
$ws.Range("L42").Value = 30984.75
$ws.Range("N42").Value = -32052.75
$ws.Range("H68").Value = 244.9
$ws.Range("J68").Value = 242.71428
$ws.Range("L68").Value = 728.14284
$ws.Range("N68").Value = -2350.14284
$ws.Range("H71").Value = 244.9
$ws.Range("J71").Value = 242.71428
$ws.Range("L71").Value = 2184.42852
$ws.Range("N71").Value = -10296.42852
$ws.Range("H140").Value = 16669752
$ws.Range("I140").Value = 18750912
$ws.Range("K140").Value = 56252736
$ws.Range("M140").Value = -56247556

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7547.2104
$ws.Range("J7").Value = 7838.8
$ws.Range("L7").Value = 7838.8
$ws.Range("N7").Value = -8062.8
$ws.Range("H126").Value = 7547.2104
$ws.Range("J126").Value = 7838.8
$ws.Range("L126").Value = 23516.4
$ws.Range("N126").Value = -28456.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 171976.27
$ws.Range("I132").Value = 2455.0195
$ws.Range("J132").Value = 1252674.2
$ws.Range("K132").Value = 7365.058499999999
$ws.Range("L132").Value = 3758022.6
$ws.Range("M132").Value = -4835.058499999999
$ws.Range("N132").Value = -3763082.6
